# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Ají" at row 860, shifting the
# existing rows 860:902 down to 861:903.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 860..902 down by one row (dimension grows from R902 to R903).
$ws.Rows("860:860").Insert()

# Populate the newly-inserted row 860 with the new record.
$ws.Range("A860").Value = 6
$ws.Range("B860").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C860").Value = 'Metropolitana'
$ws.Range("D860").Value = 44753
$ws.Range("E860").Value = 13
$ws.Range("F860").Value = 100112021
$ws.Range("G860").Value = 'Ají'
$ws.Range("H860").Value = 'Inferno'
$ws.Range("I860").Value = 'Primera'
$ws.Range("J860").Value = 70
$ws.Range("K860").Value = 10000
$ws.Range("L860").Value = 12000
$ws.Range("M860").Value = 10857
$ws.Range("N860").Value = '$/caja 12 kilos'
$ws.Range("O860").Value = 'Región de Arica y Parinacota'
$ws.Range("P860").Value = 905
$ws.Range("Q860").Value = 12
$ws.Range("R860").Value = 'Hortaliza'
